$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 801.0625
$ws.Range("J19").Value = 972.3333
$ws.Range("L19").Value = 972.3333
$ws.Range("N19").Value = -1322.3333

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 45837748
$ws.Range("J70").Value = 58828468
$ws.Range("L70").Value = 176485404
$ws.Range("N70").Value = -176485944

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 45837748
$ws.Range("J73").Value = 58828468
$ws.Range("L73").Value = 176485404
$ws.Range("N73").Value = -176487276

# Sheet ALC, row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3443.389
$ws.Range("I88").Value = 626.6667
$ws.Range("J88").Value = 4006.7334
$ws.Range("K88").Value = 626.6667
$ws.Range("L88").Value = 4006.7334
$ws.Range("M88").Value = -220.6667
$ws.Range("N88").Value = -4818.7334

# Sheet ALC, row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3443.389
$ws.Range("I91").Value = 626.6667
$ws.Range("J91").Value = 4006.7334
$ws.Range("K91").Value = 626.6667
$ws.Range("L91").Value = 4006.7334
$ws.Range("M91").Value = 777.3333
$ws.Range("N91").Value = -6814.7334

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4030.4167
$ws.Range("I100").Value = 2118.5715
$ws.Range("J100").Value = 6707
$ws.Range("K100").Value = 2118.5715
$ws.Range("L100").Value = 6707
$ws.Range("M100").Value = -1577.5715
$ws.Range("N100").Value = -7789

# Sheet ALC, row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 988.1539
$ws.Range("I103").Value = 846.25
$ws.Range("J103").Value = 1215.2
$ws.Range("K103").Value = 2538.75
$ws.Range("L103").Value = 3645.6
$ws.Range("M103").Value = -1952.75
$ws.Range("N103").Value = -4817.6

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7653.75
$ws.Range("I137").Value = 8881.308000000001
$ws.Range("K137").Value = 26643.924
$ws.Range("M137").Value = -24093.924

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3364.2727
$ws.Range("I88").Value = 2168.6667
$ws.Range("J88").Value = 3812.625
$ws.Range("K88").Value = 2168.6667
$ws.Range("L88").Value = 3812.625
$ws.Range("M88").Value = -1762.6667
$ws.Range("N88").Value = -4624.625

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3364.2727
$ws.Range("I91").Value = 2168.6667
$ws.Range("J91").Value = 3812.625
$ws.Range("K91").Value = 2168.6667
$ws.Range("L91").Value = 3812.625
$ws.Range("M91").Value = -764.6667000000002
$ws.Range("N91").Value = -6620.625

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2861.111
$ws.Range("I102").Value = 1900
$ws.Range("K102").Value = 1900
$ws.Range("M102").Value = -278

# Sheet ARM, row 111
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 60000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 60000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 60000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -68180

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4632165
$ws.Range("I122").Value = 5850364.5
$ws.Range("K122").Value = 17551093.5
$ws.Range("M122").Value = -17548643.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2010.697
$ws.Range("I132").Value = 1844.4667
$ws.Range("J132").Value = 3673
$ws.Range("K132").Value = 5533.4001
$ws.Range("L132").Value = 11019
$ws.Range("M132").Value = -3003.4001
$ws.Range("N132").Value = -16079

# Sheet ARM, row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 69664.836
$ws.Range("J133").Value = 69664.836
$ws.Range("L133").Value = 69664.836
$ws.Range("N133").Value = -74724.836

# Sheet BSM, row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 27427.715
$ws.Range("J81").Value = 27427.715
$ws.Range("L81").Value = 27427.715
$ws.Range("N81").Value = -29549.715

# Sheet BSM, row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 27427.715
$ws.Range("J84").Value = 27427.715
$ws.Range("L84").Value = 82283.145
$ws.Range("N84").Value = -92891.145

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29414542
$ws.Range("I86").Value = 2633.077
$ws.Range("K86").Value = 2633.077
$ws.Range("M86").Value = -1510.077

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 29414542
$ws.Range("I89").Value = 2633.077
$ws.Range("K89").Value = 13165.385
$ws.Range("M89").Value = -7549.385000000002

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1833.6875
$ws.Range("J94").Value = 2332.2307
$ws.Range("L94").Value = 2332.2307
$ws.Range("N94").Value = -3234.2307

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3950.0278
$ws.Range("I107").Value = 1660.9565
$ws.Range("J107").Value = 7999.923
$ws.Range("K107").Value = 1660.9565
$ws.Range("L107").Value = 7999.923
$ws.Range("M107").Value = 259.0435
$ws.Range("N107").Value = -11839.923

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3368.3
$ws.Range("I134").Value = 2866.6875
$ws.Range("K134").Value = 8600.0625
$ws.Range("M134").Value = -6065.0625

# Sheet BSM, row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

# Sheet CUL, row 48
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 995
$ws.Range("J48").Value = 1000
$ws.Range("L48").Value = 3000
$ws.Range("N48").Value = -3500

# Sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.26667
$ws.Range("I2").Value = 42.6
$ws.Range("J2").Value = 218.6
$ws.Range("K2").Value = 42.6
$ws.Range("L2").Value = 218.6
$ws.Range("M2").Value = 70.40000000000001
$ws.Range("N2").Value = -444.6

# Sheet GSM, row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 14999
$ws.Range("J40").Value = 14999
$ws.Range("L40").Value = 14999
$ws.Range("N40").Value = -15301

# Sheet GSM, row 41
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 999.5
$ws.Range("I41").Value = 999.5
$ws.Range("K41").Value = 999.5
$ws.Range("M41").Value = -644.5

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12133.667
$ws.Range("I80").Value = 21601
$ws.Range("K80").Value = 21601
$ws.Range("M80").Value = -20603

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 12133.667
$ws.Range("I83").Value = 21601
$ws.Range("K83").Value = 108005
$ws.Range("M83").Value = -103013

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6702.607
$ws.Range("I122").Value = 6755.391
$ws.Range("J122").Value = 6459.8
$ws.Range("K122").Value = 20266.173
$ws.Range("L122").Value = 19379.4
$ws.Range("M122").Value = -17816.173
$ws.Range("N122").Value = -24279.4

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2907.5
$ws.Range("I126").Value = 2761.25
$ws.Range("K126").Value = 8283.75
$ws.Range("M126").Value = -5813.75

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5549.75
$ws.Range("I82").Value = 3724.5
$ws.Range("J82").Value = 7375
$ws.Range("K82").Value = 3724.5
$ws.Range("L82").Value = 7375
$ws.Range("M82").Value = -3363.5
$ws.Range("N82").Value = -8097

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5549.75
$ws.Range("I85").Value = 3724.5
$ws.Range("J85").Value = 7375
$ws.Range("K85").Value = 3724.5
$ws.Range("L85").Value = 7375
$ws.Range("M85").Value = -2476.5
$ws.Range("N85").Value = -9871

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3396.4167
$ws.Range("I93").Value = 1829.15
$ws.Range("J93").Value = 5355.5
$ws.Range("K93").Value = 1829.15
$ws.Range("L93").Value = 5355.5
$ws.Range("M93").Value = -581.1500000000001
$ws.Range("N93").Value = -7851.5

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1007.63635
$ws.Range("I81").Value = 654.8570999999999
$ws.Range("J81").Value = 1625
$ws.Range("K81").Value = 1309.7142
$ws.Range("L81").Value = 3250
$ws.Range("M81").Value = -248.7141999999999
$ws.Range("N81").Value = -5372

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1007.63635
$ws.Range("I84").Value = 654.8570999999999
$ws.Range("J84").Value = 1625
$ws.Range("K84").Value = 6548.571
$ws.Range("L84").Value = 16250
$ws.Range("M84").Value = -1244.571
$ws.Range("N84").Value = -26858

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 503.92856
$ws.Range("I100").Value = 319.55554
$ws.Range("K100").Value = 639.11108
$ws.Range("M100").Value = -98.11108000000002

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5640.8237
$ws.Range("I122").Value = 1982.6666
$ws.Range("K122").Value = 5947.9998
$ws.Range("M122").Value = -3497.9998
